$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id changed
$ws.Range("A2").Value = 62344048

# I2, J2: cleared (Antal / Enhet removed)
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
# K2, L2: newly present but blank (Alder-Stadium / Kon)
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# S2: Noggrannhet 25 -> 10
$ws.Range("S2").Value = 10

# X2: Externid removed entirely
$ws.Range("X2").Value = ""

# Y2 / AA2: Startdatum / Slutdatum text values changed; force text so Excel
# doesn't reinterpret the ISO date string as a date serial number.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2016-08-15"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2016-08-15"

# AC2: Publik kommentar text updated
$ws.Range("AC2").Value = "övergivet finntorp med f.d. ängsmark, V om liten lada, troligen utgången, ohävdat sedan länge, senblommande"

# AD2: Ej återfunnen flag flipped to True
$ws.Range("AD2").Value = $true

# AI2: Biotop-beskrivning newly filled in
$ws.Range("AI2").Value = "ohävdad ängsmark"

# AW2 / AX2: Rapportör / Observatörer updated
$ws.Range("AW2").Value = "Tomas Troschke"
$ws.Range("AX2").Value = "Magnus Andersson"

# AY2: Projektnamn cleared
$ws.Range("AY2").Value = ""
